# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.364.61"
$ws.Range("E2").Value = "'  -3.04%  "
$ws.Range("D3").Value = "'3.731.81"
$ws.Range("E3").Value = "'  -3.85%  "
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'616.06"
$ws.Range("E5").Value = "'  +1.79%  "
$ws.Range("D6").Value = "'183.37"
$ws.Range("E6").Value = "'  +4.85%  "
$ws.Range("D7").Value = "'3.729.68"
$ws.Range("E7").Value = "'  -3.78%  "
$ws.Range("D8").Value = "'0.635"
$ws.Range("E8").Value = "'  -5.23%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("D10").Value = "'0.726"
$ws.Range("E10").Value = "'  -3.40%  "
$ws.Range("E11").Value = "'  -8.60%  "
$ws.Range("D12").Value = "'58.02"
$ws.Range("E12").Value = "'  +7.06%  "
$ws.Range("D13").Value = "'0.0000297"
$ws.Range("E13").Value = "'  -8.20%  "
$ws.Range("D14").Value = "'10.80"
$ws.Range("E14").Value = "'  -5.79%  "
$ws.Range("D15").Value = "'4.311.93"
$ws.Range("E15").Value = "'  -4.27%  "
$ws.Range("D16").Value = "'3.724.06"
$ws.Range("E16").Value = "'  -4.12%  "
$ws.Range("D17").Value = "'19.56"
$ws.Range("E17").Value = "'  -6.76%  "
$ws.Range("D18").Value = "'13.04"
$ws.Range("E18").Value = "'  -6.60%  "
$ws.Range("E19").Value = "'  -2.01%  "
$ws.Range("E20").Value = "'  -6.53%  "
$ws.Range("D21").Value = "'69.040.12"
$ws.Range("E21").Value = "'  -3.26%  "
$ws.Range("D22").Value = "'417.22"
$ws.Range("E22").Value = "'  -5.25%  "
$ws.Range("D23").Value = "'4.76"
$ws.Range("E23").Value = "'  -0.49%  "
$ws.Range("D24").Value = "'89.65"
$ws.Range("E24").Value = "'  -4.85%  "
$ws.Range("D25").Value = "'3.07"
$ws.Range("E25").Value = "'  -7.77%  "
$ws.Range("D26").Value = "'12.82"
$ws.Range("E26").Value = "'  -7.88%  "
$ws.Range("D27").Value = "'11.03"
$ws.Range("E27").Value = "'  -6.41%  "
$ws.Range("E28").Value = "'  -3.56%  "
$ws.Range("E29").Value = "'  +1.66%  "
$ws.Range("D30").Value = "'9.67"
$ws.Range("E30").Value = "'  -8.16%  "
$ws.Range("D31").Value = "'33.26"
$ws.Range("E31").Value = "'  -5.48%  "
$ws.Range("D32").Value = "'7.41"
$ws.Range("E32").Value = "'  -15.81%  "
$ws.Range("D33").Value = "'12.59"
$ws.Range("E33").Value = "'  -7.47%  "
$ws.Range("E34").Value = "'  -5.46%  "
$ws.Range("D35").Value = "'66.08"
$ws.Range("E35").Value = "'  -5.31%  "
$ws.Range("D36").Value = "'44.21"
$ws.Range("E36").Value = "'  -7.76%  "
$ws.Range("D37").Value = "'610.73"
$ws.Range("E37").Value = "'  -3.52%  "
$ws.Range("D38").Value = "'0.0₃0899"
$ws.Range("E38").Value = "'  -10.90%  "
$ws.Range("D39").Value = "'0.410"
$ws.Range("E39").Value = "'  -6.16%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.03%  "
$ws.Range("E41").Value = "'  -0.18%  "
$ws.Range("D42").Value = "'0.140"
$ws.Range("E42").Value = "'  -5.00%  "
$ws.Range("D43").Value = "'3.09"
$ws.Range("E43").Value = "'  -7.44%  "
$ws.Range("D44").Value = "'0.0445"
$ws.Range("E44").Value = "'  -5.88%  "
$ws.Range("D45").Value = "'2.69"
$ws.Range("E45").Value = "'  -6.04%  "
$ws.Range("B46").Value = "'dogwifhat"
$ws.Range("C46").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.81"
$ws.Range("E46").Value = "'  -11.15%  "
$ws.Range("B47").Value = "'THORChain"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'9.30"
$ws.Range("E47").Value = "'  -9.25%  "
$ws.Range("D48").Value = "'0.137"
$ws.Range("E48").Value = "'  -5.72%  "
$ws.Range("D49").Value = "'2.808.57"
$ws.Range("E49").Value = "'  -3.67%  "
$ws.Range("D50").Value = "'2.73"
$ws.Range("E50").Value = "'  -5.47%  "
$ws.Range("D51").Value = "'0.000269"
$ws.Range("E51").Value = "'  -3.29%  "
